$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.42
$ws.Range("G2").Value = 2.5
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 4.3
$ws.Range("J2").Value = 2.84
$ws.Range("K2").Value = 2.92
$ws.Range("L2").Value = 3.6
$ws.Range("M2").Value = 1.22
$ws.Range("N2").Value = 2.02
$ws.Range("O2").Value = 1.95
$ws.Range("P2").Value = 1.31
$ws.Range("Q2").Value = 4.1
$ws.Range("R2").Value = 1.09
$ws.Range("S2").Value = 11
$ws.Range("T2").Value = 2.86
$ws.Range("U2").Value = 1.44
$ws.Range("V2").Value = 1.32
$ws.Range("W2").Value = 1.68
$ws.Range("X2").Value = 5.9
$ws.Range("Y2").Value = 8.800000000000001
$ws.Range("Z2").Value = 26
$ws.Range("AB2").Value = 6.2
$ws.Range("AC2").Value = 7.6
$ws.Range("AD2").Value = 23
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 13.5
$ws.Range("AG2").Value = 16
$ws.Range("AH2").Value = 940
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 44
$ws.Range("AK2").Value = 60
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
